$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" (column D) values in this sheet are digit-grouped/decimal
# looking strings (e.g. "4.90", "204.36") that Excel would otherwise
# auto-convert to a number (dropping the trailing zero / reformatting),
# so force those specific cells to Text format before writing, then
# restore the original "Normal" cell style once the literal text is in
# place (keeps formatting identical to the rest of the sheet).
$textCells = @('D5', 'D6', 'D13', 'D16', 'D19', 'D20', 'D21', 'D23', 'D24', 'D25', 'D27', 'D32', 'D33', 'D35', 'D36', 'D37', 'D38', 'D39', 'D41', 'D42', 'D46', 'D48', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '76.191.53'
$ws.Range('E2').Value = '  +1.72%  '
$ws.Range('D3').Value = '2.922.74'
$ws.Range('E3').Value = '  +3.55%  '
$ws.Range('D5').Value = '204.36'
$ws.Range('E5').Value = '  +8.84%  '
$ws.Range('D6').Value = '598.21'
$ws.Range('E6').Value = '  +0.47%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  +0.19%  '
$ws.Range('E9').Value = '  +2.50%  '
$ws.Range('D10').Value = '2.921.29'
$ws.Range('E10').Value = '  +3.57%  '
$ws.Range('E11').Value = '  +16.46%  '
$ws.Range('E12').Value = '  +0.62%  '
$ws.Range('D13').Value = '4.90'
$ws.Range('E13').Value = '  +0.05%  '
$ws.Range('D14').Value = '3.456.69'
$ws.Range('E14').Value = '  +3.47%  '
$ws.Range('D15').Value = '76.062.51'
$ws.Range('E15').Value = '  +1.68%  '
$ws.Range('D16').Value = '28.08'
$ws.Range('E16').Value = '  +4.71%  '
$ws.Range('E17').Value = '  +1.37%  '
$ws.Range('D18').Value = '2.920.56'
$ws.Range('E18').Value = '  +3.55%  '
$ws.Range('D19').Value = '12.92'
$ws.Range('E19').Value = '  +4.86%  '
$ws.Range('D20').Value = '8.75'
$ws.Range('E20').Value = '  -2.13%  '
$ws.Range('D21').Value = '372.12'
$ws.Range('E21').Value = '  -1.45%  '
$ws.Range('E22').Value = '  +1.99%  '
$ws.Range('D23').Value = '4.28'
$ws.Range('E23').Value = '  +5.34%  '
$ws.Range('D24').Value = '71.56'
$ws.Range('E24').Value = '  +0.91%  '
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('D26').Value = '3.073.42'
$ws.Range('E26').Value = '  +3.65%  '
$ws.Range('D27').Value = '4.22'
$ws.Range('E27').Value = '  +1.32%  '
$ws.Range('E28').Value = '  -0.24%  '
$ws.Range('E29').Value = '  +3.49%  '
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('E31').Value = '  +1.48%  '
$ws.Range('D32').Value = '501.53'
$ws.Range('E32').Value = '  -3.27%  '
$ws.Range('D33').Value = '7.76'
$ws.Range('E33').Value = '  +0.28%  '
$ws.Range('E34').Value = '  +2.47%  '
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').Value = '165.16'
$ws.Range('E36').Value = '  +1.43%  '
$ws.Range('D37').Value = '20.25'
$ws.Range('E37').Value = '  +1.55%  '
$ws.Range('B38').Value = 'Cronos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D38').Value = '0.107'
$ws.Range('E38').Value = '  +25.11%  '
$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').Value = '19.62'
$ws.Range('E39').Value = '  +1.33%  '
$ws.Range('E40').Value = '  -4.77%  '
$ws.Range('D41').Value = '0.364'
$ws.Range('E41').Value = '  +6.68%  '
$ws.Range('D42').Value = '182.74'
$ws.Range('E42').Value = '  -2.13%  '
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('E45').Value = '  -0.65%  '
$ws.Range('D46').Value = '40.00'
$ws.Range('E47').Value = '  -1.42%  '
$ws.Range('D48').Value = '2.37'
$ws.Range('E48').Value = '  +1.85%  '
$ws.Range('E49').Value = '  -1.13%  '
$ws.Range('E50').Value = '  +0.12%  '
$ws.Range('D51').Value = '22.45'
$ws.Range('E51').Value = '  +7.30%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
